$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix the JSON-import bug: column B ("Title") was always populated with the
# literal header text "Title" instead of each book's actual title.
$ws.Range("B2").Value2 = "Everyday Italian"
$ws.Range("B3").Value2 = "Harry Potter"
$ws.Range("B4").Value2 = "XQuery Kick Start"
$ws.Range("B5").Value2 = "Learning XML"

# Widen column B so the now-longer titles are readable.
$ws.Columns.Item(2).ColumnWidth = 16.4
